$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Coin name / link swaps (rows 11-12 and 28-31 re-sorted) ---
$ws.Range('B11').Value = 'Toncoin'
$ws.Range('C11').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('B12').Value = 'TRON'
$ws.Range('C12').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('B28').Value = 'Fetch.AI'
$ws.Range('C28').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('B29').Value = 'Binance-PegBSC-USD'
$ws.Range('C29').Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range('B30').Value = 'InternetComputer(DFINITY)'
$ws.Range('C30').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('B31').Value = 'Bittensor'
$ws.Range('C31').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'

# --- Price (column D) and Volume(1h) (column E) updates ---
# Force text format on touched D/E cells so numeric-looking strings
# (e.g. '0.606', '1.00') are preserved as text, matching the source file's
# inlineStr cell type instead of being auto-converted to numbers.
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '61.506.00'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  -3.05%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.478.53'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  -5.40%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.00'
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  +0.02%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '553.64'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  -3.44%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '146.11'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  -5.15%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '1.00'
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  -0.04%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.606'
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  -2.54%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '2.480.64'
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  -5.21%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.108'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  -7.11%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '5.45'
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  -5.13%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.154'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  -1.45%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.357'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  -5.37%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '26.21'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  -6.46%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '2.925.66'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  -5.39%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.0000167'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  -8.02%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '61.424.34'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  -3.11%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '2.471.19'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '11.17'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  -6.77%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '6.97'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  -7.26%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '4.23'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  -5.97%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '322.19'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  -5.43%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '1.00'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  +0.01%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '1.81'
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  -2.31%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '64.08'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  -5.27%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.0₃0989'
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  -7.27%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '2.591.48'
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  -5.88%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '1.51'
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  -3.31%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.999'
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  -0.15%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '8.44'
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  -7.00%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '529.68'
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  -8.39%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '7.66'
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  -2.23%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.152'
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  -5.33%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.91'
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  -6.10%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.58'
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  -8.28%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '5.91'
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  -9.84%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '4.91'
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  -7.77%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.00'
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  +0.08%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.385'
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  -3.26%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '18.57'
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  -5.40%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '147.99'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  -2.76%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '1.74'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  -7.70%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '40.32'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  -3.02%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '2.34'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  -5.89%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '148.40'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  -5.80%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '3.60'
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  -6.81%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '20.75'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  -12.47%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.0534'
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  -7.95%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.597'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  -4.68%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.0948'
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  -4.48%  '
